$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Artesanal Investimentos'
$ws.Cells.Item(2, 2).Value = 'Estágio de Risco - Mercado Financeiro'
$ws.Cells.Item(2, 3).Value = 'Presencial'
$ws.Cells.Item(2, 4).Value = 'Microsoft Excel, Conhecimentos em Lógica de Programação, Cursando os últimos 4 semestres de engenharias, matemática, estatística, economia ou física'
$ws.Cells.Item(2, 5).Value = 'Plataforma Educacional, Clínica Corporativa, VT sem desconto ou estacionamento, Gympass, Assistência Médica, Vale Refeição, Cartão Flash, Sala de Leitura, Auxílio Creche, Empréstimo Consignado, Freshbook, Bônus'
$ws.Cells.Item(2, 6).Value = 'Estudantes a partir do 5º semestre de engenharias, matemática, estatística, economia ou física'
$ws.Cells.Item(2, 7).Value = 'Mercado Financeiro, Risco'
$ws.Cells.Item(2, 8).Value = 'Processamento do cálculo de risco das carteiras dos fundos, Análise de risco, Atualização de apresentações para o Comitê Risco, Monitoramento de garantias e resseguro.'

$ws.Cells.Item(3, 1).Value = 'Empresa de consultoria para o setor de saúde'
$ws.Cells.Item(3, 2).Value = 'Estágio'
$ws.Cells.Item(3, 3).Value = 'Não mencionado'
$ws.Cells.Item(3, 4).Value = 'Formação em Engenharia, Administração ou Economia; Conhecimento do Pacote Office; Boa capacidade de comunicação; Perfil focado, dinâmico e com capacidade analítica'
$ws.Cells.Item(3, 5).Value = 'Bolsa competitiva com política de bônus agressiva'
$ws.Cells.Item(3, 6).Value = 'Não mencionado'
$ws.Cells.Item(3, 7).Value = 'Consultoria, saúde, finanças'
$ws.Cells.Item(3, 8).Value = 'Modelagem financeira, estruturação e otimização de processos, interação direta com os clientes.'

$ws.Cells.Item(4, 1).Value = 'Finlead'
$ws.Cells.Item(4, 2).Value = 'Estágio em Análise de Investimentos e Controle de Operações'
$ws.Cells.Item(4, 3).Value = 'Não especificado'
$ws.Cells.Item(4, 4).Value = 'Proatividade, perfil analítico e quantitativo, conhecimento avançado em Excel e matemática financeira, graduação em Engenharia ou finanças (preferencialmente administração, ciências contábeis ou economia), conhecimento em finanças, inglês avançado, conhecimento em controladoria de fundos e regulamentação aplicável à gestão de fundos de investimento, domínio do pacote Office'
$ws.Cells.Item(4, 5).Value = 'Bolsa competitiva ao mercado, Vale Refeição, Vale Transporte'
$ws.Cells.Item(4, 6).Value = 'Estudantes universitários em Engenharia ou Finanças'
$ws.Cells.Item(4, 7).Value = 'Análise de Investimentos e Controle de Operações'
$ws.Cells.Item(4, 8).Value = 'Elaboração de relatórios, acompanhamento de aquisições, monitoramento de aderência de compras, preparação de relatórios gerenciais, atuação no relacionamento com empresas, integração com a originadora do fundo.'

$ws.Cells.Item(5, 1).Value = 'Jera Capital Family Office'
$ws.Cells.Item(5, 2).Value = 'Estágio'
$ws.Cells.Item(5, 3).Value = 'Não especificada'
$ws.Cells.Item(5, 4).Value = 'Inglês Avançado ou Fluente, Graduação em andamento nas áreas de Economia, Administração, Engenharia ou áreas Correlatas, Python Intermediário, Pacote Office Avançado (VBA), Disponibilidade para modelo de trabalho (3x presencial na semana).'
$ws.Cells.Item(5, 5).Value = 'Não especificada'
$ws.Cells.Item(5, 6).Value = 'Estudantes nas áreas de Economia, Administração, Engenharia ou áreas Correlatas'
$ws.Cells.Item(5, 7).Value = 'Gestão de investimentos, Automatização de processos e rotinas diárias'
$ws.Cells.Item(5, 8).Value = 'Auxílio ao time, buscar automatizar os processos e rotinas do dia a dia, atuação 360° dentro da gestora.'

$ws.Cells.Item(6, 1).Value = 'Fram Capital'
$ws.Cells.Item(6, 2).Value = 'Estagiário Financeiro'
$ws.Cells.Item(6, 3).Value = 'Vila Nova Conceição – Zona Sul'
$ws.Cells.Item(6, 4).Value = 'Cursando Engenharias, a partir do 5º semestre, entendimento acerca de contabilidade (balanço, DRE, etc)'
$ws.Cells.Item(6, 5).Value = 'Bolsa Auxílio R$2.600,00, Vale Refeição R$600,00'
$ws.Cells.Item(6, 6).Value = 'Interessados enviar CV com o assunto “Financeiro” para jcosta@framcapital.com'
$ws.Cells.Item(6, 7).Value = 'Financeira'
$ws.Cells.Item(6, 8).Value = 'Rotinas Administrativas; Emissão de Notas Fiscais; Contas a pagar; Contas a receber; Atualização de sistema financeiro; Conciliação Bancária; Automatização de Processos.'

$ws.Cells.Item(7, 1).Value = 'Paggo (Stealth Startup)'
$ws.Cells.Item(7, 2).Value = 'Estagiário em Engenharia de Software'
$ws.Cells.Item(7, 3).Value = 'Remoto'
$ws.Cells.Item(7, 4).Value = 'Adquirir conhecimento e proficiência no trabalho com as stacks de tecnologia mais modernas do mercado, desenvolver produtos de software usando melhores práticas de lógica e programação, influenciar o desenvolvimento de produto.'
$ws.Cells.Item(7, 5).Value = 'R$ 3.900,00 (R$ 3000,00 contrato de estágio + R$ 900,00 em um cartão de benefícios flexíveis)'
$ws.Cells.Item(7, 6).Value = 'Jovens brilhantes e ambiciosos'
$ws.Cells.Item(7, 7).Value = 'Engenharia de Software'
$ws.Cells.Item(7, 8).Value = 'Adquirir conhecimento e proficiência no trabalho com as stacks de tecnologia mais modernas do mercado, desenvolver produtos de software usando melhores práticas de lógica e programação, influenciar o desenvolvimento de produto.'

$ws.Cells.Item(8, 1).Value = 'Grupo Boticário'
$ws.Cells.Item(8, 2).Value = 'Pessoa Desenvolvedora BackEnd Java/Kotlin Especialista II (Engenharia de Crédito) (Produtos Digitais Financeiros)'
$ws.Cells.Item(8, 3).Value = 'Remoto'
$ws.Cells.Item(8, 4).Value = 'Experiência com soluções cloud AWS, conhecimento em Python, JavaScript/TypeScript(Node) e/ou Java/Kotlin, interesse e conhecimento em estruturas de dados, experiência com monitoração e logging, conhecimento de Rest e orientação a eventos(Kafka), experiência com testes unitários e de integração, visão de qualidade de software, domínio do git ou outro sistema de controle de versão colaborativo, conhecimento em CI / CD, conhecimento em serverless framework, familiaridade em banco de dados relacional e não relacional, conhecimento em produtos financeiros voltados à créditos, capacidade de conduzir root cause analysis em problemas de software, experiência com testes de carga.'
$ws.Cells.Item(8, 5).Value = 'A combinar'
$ws.Cells.Item(8, 6).Value = 'Vagas destinadas aos grupos minorizados priorizados em nossa estratégia: pessoas com deficiência, pessoas negras (pretas e pardas), mulheres (cis e trans), pessoas da comunidade LGBTQIA+ e pessoas 50+.'
$ws.Cells.Item(8, 7).Value = 'Desenvolvimento de Plataforma de Crédito, Produtos Digitais Financeiros.'
$ws.Cells.Item(8, 8).Value = 'Participar do desenho da arquitetura de solução e cenários de uso, definir requisitos técnicos, arquitetura de sistemas e melhores abordagens de desenvolvimento, apoiar profissionais que precisem da sua experiência ou orientação, contribuir com os time de infraestrutura, segurança e arquitetura para determinar as melhores soluções para os problemas, fornecer orientação técnica e suporte à equipe, compartilhar conhecimento e melhores práticas, ajudar na evolução técnica dos outros desenvolvedores.'

$ws.Cells.Item(9, 1).Value = 'Cadastra'
$ws.Cells.Item(9, 2).Value = 'Data Strategy Assistant - Estágio'
$ws.Cells.Item(9, 3).Value = 'Remoto'
$ws.Cells.Item(9, 4).Value = 'Habilidades de comunicação interpessoal, interesse em aprender e se desafiar, desejo de trabalhar em um ambiente dinâmico, matrícula ativa em um curso superior (tecnologia, administração, publicidade, marketing, engenharias ou áreas correlatas). Diferenciais: Conhecimentos básicos em Ferramentas de Analytics, Data Visualization, HTML e linguagens de programação (Javascript), conhecimento de métricas de marketing digital, GTM e tagueamento, inglês para leitura.'
$ws.Cells.Item(9, 5).Value = 'A combinar'
$ws.Cells.Item(9, 6).Value = 'Estudantes com matrícula ativa em curso superior'
$ws.Cells.Item(9, 7).Value = 'Data & Analytics'
$ws.Cells.Item(9, 8).Value = 'Participar de reuniões com clientes e equipes internas, responsável pela implementação, gestão e manutenção de tags e plataformas de Digital Analytics para Sites e Aplicativos, responsável pela criação e manutenção de bases de dados, responsável por criar documentos instrutivos de implementações para clientes, responsável por planejar e gerenciar o plano de métricas, auxiliar na confecção de dashboards e relatórios básicos para clientes, auxiliar no desenvolvimento de análises descritivas e diagnósticas para os clientes.'

$ws.Cells.Item(10, 1).Value = 'Santander'
$ws.Cells.Item(10, 2).Value = 'Programa de Estágio Santander'
$ws.Cells.Item(10, 3).Value = 'Remoto'
$ws.Cells.Item(10, 4).Value = 'Ser estudante de graduação ou tecnólogo a partir do 2° semestre, disponibilidade para jornada de 4h ou 6h diárias, ser curioso, questionador e com vontade de transformar.'
$ws.Cells.Item(10, 5).Value = 'A combinar'
$ws.Cells.Item(10, 6).Value = 'Estudantes de graduação ou tecnólogos a partir do 2° semestre'
$ws.Cells.Item(10, 7).Value = 'Lojas, Corporativo, áreas de apoio (Tecnologia, Comunicação, Riscos, Jurídico, entre outras), Atacado.'
$ws.Cells.Item(10, 8).Value = 'Participação em projetos importantes para o banco, desenvolvimento de potencial em projetos que têm o poder de transformar o mercado, a sociedade e a vida de nossos clientes.'

$ws.Cells.Item(11, 1).Value = 'Akross'
$ws.Cells.Item(11, 2).Value = 'Estagiária em Desenvolvimento Backend'
$ws.Cells.Item(11, 3).Value = 'Remoto'
$ws.Cells.Item(11, 4).Value = 'Cursando ensino superior em Ciência da Computação, Engenharia de Software, Sistemas de Informação, Análise e Desenvolvimento de Sistemas ou áreas correlatas; Formatura prevista a partir de Junho/2026; Conhecimento em Orientação a Objetos; Conhecimento em Java 8+ e Spring Boot; Conhecimento em Hibernate/JPA; Conhecimento em bancos de dados relacionais (MySQL, PostgreSQL) ou NoSQL (MongoDB); Capacidade de atuar em equipe e boa comunicação.'
$ws.Cells.Item(11, 5).Value = 'A combinar'
$ws.Cells.Item(11, 6).Value = 'Estudantes de Ciência da Computação, Engenharia de Software, Sistemas de Informação, Análise e Desenvolvimento de Sistemas ou áreas correlatas.'
$ws.Cells.Item(11, 7).Value = 'Desenvolvimento Backend'
$ws.Cells.Item(11, 8).Value = 'Participar de forma supervisionada no desenvolvimento de aplicações Java, sob arquitetura de microsserviços; Colaborar com a equipe de desenvolvimento em algumas fases do ciclo de vida do software; Contribuir com a documentação técnica dos sistemas e aplicações; Estar junto do time na resolução de problemas das aplicações para absorção de conhecimentos; Trabalhar com a metodologia ágil Scrum, participando ativamente das cerimônias.'

$ws.Cells.Item(12, 1).Value = 'Honda'
$ws.Cells.Item(12, 2).Value = 'Estágio TI - Desenvolvimento'
$ws.Cells.Item(12, 3).Value = 'Remoto (presencial na Honda Morumbi 1x por mês), São Paulo'
$ws.Cells.Item(12, 4).Value = 'Cursando graduação em tecnologia da informação e correlatas; Conhecimento em Informática: Linguagens de programação: Java, Genexus, web banco de dados: DB2(AS400 e/ou AIX) e Cloud; Inglês intermediário.'
$ws.Cells.Item(12, 5).Value = 'A combinar'
$ws.Cells.Item(12, 6).Value = 'Estudantes de graduação em tecnologia da informação e correlatas'
$ws.Cells.Item(12, 7).Value = 'TI, Desenvolvimento'
$ws.Cells.Item(12, 8).Value = 'Apoio na organização da entrega contínua dos produtos de software ou projetos de sistemas; Apoio aos times de desenvolvedores na aplicação das melhorias práticas e técnicas de codificação; Apoio na gestão dos fornecedores externos.'

$ws.Cells.Item(13, 1).Value = 'Radix Engenharia e Software'
$ws.Cells.Item(13, 2).Value = 'Estágio em Desenvolvimento de Software'
$ws.Cells.Item(13, 3).Value = 'Remoto'
$ws.Cells.Item(13, 4).Value = 'Cursando graduação em Ciência da Computação, Engenharia de Software, Sistemas de Informação ou áreas afins com previsão conclusão para 2026/1, conhecimento básico de programação em Python, familiaridade com JavaScript e frameworks frontend como React.js, desejo de aprender sobre bancos de dados NoSQL, especialmente MongoDB.'
$ws.Cells.Item(13, 5).Value = 'A combinar'
$ws.Cells.Item(13, 6).Value = 'Estudantes de graduação em Ciência da Computação, Engenharia de Software, Sistemas de Informação ou áreas afins'
$ws.Cells.Item(13, 7).Value = 'Desenvolvimento de Software, Engenharia de Software, Sistemas de Informação'
$ws.Cells.Item(13, 8).Value = 'Auxiliar no desenvolvimento e manutenção da infraestrutura backend utilizando Python e frameworks como Flask ou FastAPI, suportar a implementação e gerenciamento de bancos de dados MongoDB, participar do desenvolvimento de interfaces de usuário utilizando JavaScript, React.js, Node.js e Express, ajudar na utilização de Redux para gerenciamento de estado em aplicações React, colaborar com a equipe de desenvolvimento para resolver problemas e implementar melhorias, participar de reuniões e sessões de brainstorming para contribuir com ideias inovadoras, manter a documentação técnica organizada e atualizada.'

$ws.Cells.Item(14, 1).Value = 'Bradesco'
$ws.Cells.Item(14, 2).Value = 'Programa de Estágio Bradesco 2024 Atacado'
$ws.Cells.Item(14, 3).Value = 'Remoto'
$ws.Cells.Item(14, 4).Value = 'Estudante de nível superior (bacharelado, licenciatura ou tecnólogo), cursando a partir do 2º semestre da graduação ou do 1º semestre de tecnólogo; Disponibilidade para uma jornada de 20, 25 ou 30 horas semanais e, preferencialmente, com possibilidade de estagiar por 2 anos; Inglês a partir do nível intermediário'
$ws.Cells.Item(14, 5).Value = 'A combinar'
$ws.Cells.Item(14, 6).Value = 'Estudantes de nível superior a partir do 2º semestre'
$ws.Cells.Item(14, 7).Value = 'Área financeira, Atacado'
$ws.Cells.Item(14, 8).Value = 'Mergulhar no universo dos grandes investidores, desvendando as soluções financeiras mais inovadoras do mercado, atendimento a grandes investidores institucionais do Brasil e clientes de private banking.'

$ws.Cells.Item(15, 1).Value = 'Britvic Brasil'
$ws.Cells.Item(15, 2).Value = 'Estagiária em Comércio Exterior'
$ws.Cells.Item(15, 3).Value = 'Remoto'
$ws.Cells.Item(15, 4).Value = 'Cursando Administração, Logística, Comércio Exterior, Relações Internacionais ou cursos afins; Conhecimento intermediário em Pacote Office - Word, Excel e Power Point; Conhecimento intermediário em inglês é desejável.'
$ws.Cells.Item(15, 5).Value = 'A combinar; Bolsa-Auxílio; Vale Transporte; Vale Refeição; Wellhub (gympass)'
$ws.Cells.Item(15, 6).Value = 'Estudantes dos cursos mencionados nos requisitos'
$ws.Cells.Item(15, 7).Value = 'Comércio Exterior, Administração, Logística, Relações Internacionais'
$ws.Cells.Item(15, 8).Value = 'Acompanhamento do fluxo de exportação, preparação do pacote de documentos dos embarques, acompanhamento de coletas realizadas dentro do território brasileiro, responsável pelo fluxo de pagamentos dos fornecedores da logísticas.'

$ws.Cells.Item(16, 1).Value = 'alt.bank'
$ws.Cells.Item(16, 2).Value = 'Estagiária em Marketing'
$ws.Cells.Item(16, 3).Value = 'Remoto'
$ws.Cells.Item(16, 4).Value = 'Graduação em andamento em cursos correlatos a Publicação, Marketing ou Relações Públicas, habilidades de planejamento de comunicação e execução de redes sociais, aptidão para criar conteúdo envolvente para mídias sociais, experiência com algumas plataformas de mídia social relevantes (Facebook, Instagram e Tiktok), desejável experiência com Linkedin, Twitter e Pinterest, desejável ter seu próprio perfil relevante em alguma rede social, fortes habilidades de comunicação verbal e escrita em português, bom domínio da língua inglesa.'
$ws.Cells.Item(16, 5).Value = 'A combinar, inclui bolsa estágio, vale alimentação/refeição, Gympass, plano de saúde e odontológico.'
$ws.Cells.Item(16, 6).Value = 'Estudantes de graduação em cursos correlatos a Publicação, Marketing ou Relações Públicas.'
$ws.Cells.Item(16, 7).Value = 'Marketing, Mídias Sociais.'
$ws.Cells.Item(16, 8).Value = 'Auxiliar na concepção e entrega de estratégias de mídia social, criar planejamento e calendário de postagens, publicar e dar manutenção nos conteúdos publicados nas redes sociais, auxiliar no desenvolvimento, lançamento e gerenciamento de novas campanhas, gerar relatórios e analisar o desempenho em plataformas de mídia social, auxiliar na identificação de tendências de consumo, ajudar a otimizar o conteúdo para incentivar a interação e o envolvimento da comunidade, pesquisar e avaliar as ferramentas e técnicas mais recentes para melhor medir a atividade nas redes sociais.'

$ws.Cells.Item(17, 1).Value = 'ACE Ventures'
$ws.Cells.Item(17, 2).Value = 'Estagiário(a) de Consultoria - Estratégia e Inteligência'
$ws.Cells.Item(17, 3).Value = 'Remoto'
$ws.Cells.Item(17, 4).Value = 'Cursando graduação em Administração, Economia, Engenharia ou áreas relacionadas (a partir do 3º semestre); Interesse genuíno por empreendedorismo, inovação e novas tecnologias; Conhecimento em estratégia empresarial, pesquisa de mercado e análise de dados; Mente curiosa e analítica, com facilidade para identificar padrões e tendências; Excelente comunicação oral e escrita.'
$ws.Cells.Item(17, 5).Value = 'A combinar'
$ws.Cells.Item(17, 6).Value = 'Estudantes a partir do 3º semestre dos cursos de Administração, Economia, Engenharia ou áreas relacionadas.'
$ws.Cells.Item(17, 7).Value = 'Consultoria, Estratégia Corporativa, Pesquisa de Mercado, Análise de Empresas, Projetos de Consultoria.'
$ws.Cells.Item(17, 8).Value = 'Auxiliar na formulação e implementação de estratégias corporativas; Conduzir pesquisas de mercado; Elaborar apresentações executivas; Organizar e gerenciar informações do projeto; Realizar análises financeiras e estratégicas de empresas; Participar de todas as etapas dos projetos de consultoria.'

$ws.Cells.Item(18, 1).Value = 'BTG Pactual'
$ws.Cells.Item(18, 2).Value = 'Estágio Short - Automation 2024'
$ws.Cells.Item(18, 3).Value = 'Remoto'
$ws.Cells.Item(18, 4).Value = 'Cursando regularmente um curso de formação superior (nível bacharel ou tecnólogo), possuir um computador com conexão à internet, disponibilidade para trabalhar no mínimo 2 meses, disponibilidade para trabalhar presencialmente em São Paulo ou Rio de Janeiro em caso de efetivação, conhecimento em base lógica de programação.'
$ws.Cells.Item(18, 5).Value = 'A combinar'
$ws.Cells.Item(18, 6).Value = 'Estudantes de curso superior (nível bacharel ou tecnólogo)'
$ws.Cells.Item(18, 7).Value = 'Automation'
$ws.Cells.Item(18, 8).Value = 'Desenvolver automações, construir telas de front-end em ReactJS e Flask, programar back-end e APIs de serviços em Python, criar processos padronizados, lidar com produtos financeiros, desenvolver soluções definidas, prestar suporte, manter um relacionamento próximo com os usuários, negociar a prioridade das soluções, fomentar a missão da área, documentar os processos e soluções, acompanhar o time de Discovery, desenvolver soluções de alto impacto, garantir a estabilidade e o funcionamento das soluções do time.'

$ws.Cells.Item(19, 1).Value = 'Arquivei'
$ws.Cells.Item(19, 2).Value = 'Estágio em CRM (Marketing)'
$ws.Cells.Item(19, 3).Value = 'Remoto'
$ws.Cells.Item(19, 4).Value = 'Estar cursando Ensino Superior nas áreas de marketing, engenharia, administração, data science, ou similares; habilidades de comunicação oral e escrita; habilidade para trabalhar com Google Sheet; Noções de marketing digital.'
$ws.Cells.Item(19, 5).Value = 'A combinar'
$ws.Cells.Item(19, 6).Value = 'Estudantes de Ensino Superior nas áreas de marketing, engenharia, administração, data science, ou similares.'
$ws.Cells.Item(19, 7).Value = 'Marketing, CRM.'
$ws.Cells.Item(19, 8).Value = 'Auxiliar no planejamento, implementação e análise de desempenho de réguas de comunicação e campanhas; Auxiliar na criação de conteúdo das comunicações; Contribuir no desenvolvimento de relatórios de performance para identificação de oportunidades de otimização; Contribuir com o aumento da representatividade de CRM na geração de demanda para o time de vendas, e com a melhora das principais métricas de CRM.'

$ws.Cells.Item(20, 1).Value = 'Bain & Company'
$ws.Cells.Item(20, 2).Value = 'AC - Associate Consultant, ACI - Associate Consultant Intern, Summer Associate Consultant Intern'
$ws.Cells.Item(20, 3).Value = 'Remoto'
$ws.Cells.Item(20, 4).Value = 'Estudante universitário ou recém-formado, habilidades analíticas, interpessoais, criativas, de resolução de problemas e de liderança'
$ws.Cells.Item(20, 5).Value = 'A combinar'
$ws.Cells.Item(20, 6).Value = 'Estudantes universitários ou recém-formados'
$ws.Cells.Item(20, 7).Value = 'Consultoria Estratégica'
$ws.Cells.Item(20, 8).Value = 'Trabalho em equipe, identificação de fontes de informação, coleta e interpretação de dados, execução de análises, apresentação de resultados, entrevistas com consumidores, concorrentes, fornecedores e empregadores, supervisão de colegas mais novos.'

$ws.Cells.Item(21, 1).Value = '+A Educação'
$ws.Cells.Item(21, 2).Value = 'Estágio em Talent Acquisition (Atração e Seleção)'
$ws.Cells.Item(21, 3).Value = 'Remoto'
$ws.Cells.Item(21, 4).Value = 'Ensino Superior em andamento, habilidades com o Pacote Office, disponibilidade para atuar em formato de trabalho híbrido em Porto Alegre/RS ou remota em outros locais, proatividade, comunicação, flexibilidade, organização e bom relacionamento interpessoal.'
$ws.Cells.Item(21, 5).Value = 'A combinar'
$ws.Cells.Item(21, 6).Value = 'Estudantes com ensino superior em andamento'
$ws.Cells.Item(21, 7).Value = 'Recursos Humanos, área de Desenvolvimento, Atração e Seleção'
$ws.Cells.Item(21, 8).Value = 'Mapeamento de talentos, busca ativa de candidatos, alinhamento de perfil das vagas, divulgação, triagem, entrevistas, processo de admissão de novos colaboradores, gestão das plataformas de seleção.'

$ws.Cells.Item(22, 1).Value = 'Rehagro'
$ws.Cells.Item(22, 2).Value = 'Estágio - Área Comercial'
$ws.Cells.Item(22, 3).Value = 'Remoto'
$ws.Cells.Item(22, 4).Value = 'Estudante dos cursos de Ciências Agrárias, conhecimento na área de interesse, vocação e vontade de se especializar na área comercial, disponibilidade de eventualmente estar presente no laboratório nas unidades de MG ou SP.'
$ws.Cells.Item(22, 5).Value = 'A combinar'
$ws.Cells.Item(22, 6).Value = 'Estudantes dos cursos de Ciências Agrárias, preferencialmente do 4º ao 8º período'
$ws.Cells.Item(22, 7).Value = 'Comercial, Vendas Internas, Ciências Agrárias'
$ws.Cells.Item(22, 8).Value = 'Dar apoio nos processo de vendas internas, em atividades de contato direto com clientes, acompanhamento de pedidos e fornecimento de suporte ao pós-venda, aprender sobre os produtos e serviços oferecidos pelo laboratório, entender as necessidades e demandas dos clientes do setor de pecuária e agricultura.'

$ws.Cells.Item(23, 1).Value = 'Visagio'
$ws.Cells.Item(23, 2).Value = 'Estágio Desenvolvedor(a) de Software'
$ws.Cells.Item(23, 3).Value = 'Remoto'
$ws.Cells.Item(23, 4).Value = 'Conhecimento básico em metodologias ágeis, padrões de projeto e arquitetura de sistemas; conhecimento básico em pelo menos uma linguagem de programação (ex: C#, Kotlin, Ruby, etc.) e framework (ex: .NET, SpringBoot, Rails, etc.)'
$ws.Cells.Item(23, 5).Value = 'A combinar'
$ws.Cells.Item(23, 6).Value = 'Universitários dos cursos de Ciência da Computação, Engenharias, Sistemas de Informação e afins'
$ws.Cells.Item(23, 7).Value = 'Engenharia de Dados, Automação de Processos, Gestão de TI e Desenvolvimento de Software'
$ws.Cells.Item(23, 8).Value = 'Desenvolvimento front-end e/ou back-end; Correção de bugs e melhoria contínua; Definição de melhores tecnologias a serem usadas nos sistemas desenvolvidos; Propor e implementar novas ferramentas, técnicas e metodologias; Compartilhar e evoluir o conhecimento técnico do time.'

$ws.Cells.Item(24, 1).Value = 'Smarthis'
$ws.Cells.Item(24, 2).Value = 'Programa de Estágio Smarthis 2024'
$ws.Cells.Item(24, 3).Value = 'Remoto'
$ws.Cells.Item(24, 4).Value = 'Graduação a partir do 6º período, conhecimento em alguma linguagem de programação (.NET, Python, PHP, Java, C#, VBA, VB Script, entre outras), inglês avançado.'
$ws.Cells.Item(24, 5).Value = 'A combinar'
$ws.Cells.Item(24, 6).Value = 'Estudantes de graduação a partir do 6º período com conhecimentos em programação e inglês avançado.'
$ws.Cells.Item(24, 7).Value = 'RPA (Robotic Process Automation), programação, Business Analytics.'
$ws.Cells.Item(24, 8).Value = 'Desenvolver soluções em RPA em conjunto com seu time para clientes nacionais e/ou internacionais, contribuir para uma entrega de soluções de automação, ajudando aos clientes nas suas transformações digitais.'

$ws.Cells.Item(25, 1).Value = 'Radix'
$ws.Cells.Item(25, 2).Value = 'Estágio em Desenvolvimento de Negócios'
$ws.Cells.Item(25, 3).Value = 'Remoto'
$ws.Cells.Item(25, 4).Value = 'Cursando Engenharia de Processos, Mecânica, Automação ou áreas correlatas com formação prevista para a partir de 2025/2; Inglês Avançado/Fluente; Excel Intermediário/Avançado; Capacidade de gerenciamento de tempo e resolução de problemas; Vontade de aprender e se desenvolver. Diferencial: Conhecimentos da indústria de Óleo & Gás; Conhecimentos básicos de automação.'
$ws.Cells.Item(25, 5).Value = 'A combinar'
$ws.Cells.Item(25, 6).Value = 'Estudantes de Engenharia de Processos, Mecânica, Automação ou áreas correlatas.'
$ws.Cells.Item(25, 7).Value = 'Desenvolvimento de Negócios.'
$ws.Cells.Item(25, 8).Value = 'Apoiar no mapeamento e prospecção de oportunidades; Realizar pesquisas de mercado e ações de pré-venda; Elaborar e apresentar propostas técnicas-comerciais nas áreas de atuação da empresa; Acompanhar a equipe técnica nas estimativas de recursos para projetos; Representar a empresa em visitas técnicas e reuniões para a discussão de escopo de projetos; Apoiar na construção, formalização e conclusão de propostas comerciais e orçamentárias; Elaborar materiais de divulgação de soluções desenvolvidas pela Radix.'

$ws.Cells.Item(26, 1).Value = 'Liv Up'
$ws.Cells.Item(26, 2).Value = 'Estágio em Growth - Business Intelligence'
$ws.Cells.Item(26, 3).Value = 'Remoto - Anywhere Office (qualquer lugar do Brasil)'
$ws.Cells.Item(26, 4).Value = 'Formatura a partir de Dez/25, capacidade de resolução de problemas de forma analítica e criativa, pró-atividade, boa comunicação, Excel & PowerPoint avançado, Inglês avançado, SQL e Python como diferenciais'
$ws.Cells.Item(26, 5).Value = 'A combinar'
$ws.Cells.Item(26, 6).Value = 'Estudantes com formatura a partir de Dez/25'
$ws.Cells.Item(26, 7).Value = 'Growth, Business Intelligence'
$ws.Cells.Item(26, 8).Value = 'Elaboração de modelos de projeção de receita e custos, estratégia e implementação de campanhas de marketing, análise de performance de testes de Growth, automatização de processos recorrentes, alinhamento de planos de ação com times responsáveis.'
